# Automatische test-sync: 2025-07-27 19:37:50
#
# Appends Testmail #12 ("Ik heb nog geen geld terug.") as a new row 14 on the
# "Logs" sheet, extends the conditional-formatting ranges that previously
# stopped at row 13 so they cover row 14 too, and updates the "Dashboard"
# summary table so that "Retour / Terugbetaling" (now 2 occurrences) sorts
# above "Bestelling / Levering" (still 1 occurrence).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new testmail row
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Ik heb nog geen geld terug."
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Testmail #12: Ik heb nog geen geld terug."
$logs.Range("D14").Value = "Retour / Terugbetaling"
$logs.Range("E14").Value = "Geachte klant,`nBedankt voor uw bericht. Als u nog geen geld terug heeft ontvangen, zou dat kunnen betekenen dat de terugbetaling nog niet is verwerkt. Zou u ons uw transactiegegevens en/of bestelnummer kunnen doorgeven, zodat we dit verder kunnen onderzoeken en u van dienst kunnen zijn?`nMet vriendelijke groet,`n[Naam van het bedrijf] E-mailassistent"
$logs.Range("F14").Value = "2025-07-27 19:36:55"
$logs.Range("G14").Value = "Ja"
$logs.Range("H14").Value = "Nee"
$logs.Range("I14").Value = "Ja"
$logs.Range("J14").Value = "Nee"

# Setting the multi-line "Antwoord" text auto-expands the row height; put it
# back to the sheet's implicit default (matches the other data rows, which
# carry no row-height override in the source file).
$logs.Rows.Item(14).AutoFit()

# ---------------------------------------------------------------------
# 2) Logs sheet: extend the conditional formatting ranges D/G/H/I/J from
#    row 13 down to row 14 (without disturbing priorities / dxfId links)
# ---------------------------------------------------------------------
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "13")
    $newRange = $logs.Range($col + "2:" + $col + "14")
    $fc = $oldRange.FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($newRange)
}

# ---------------------------------------------------------------------
# 3) Dashboard sheet: row 5 and row 6 swap places - "Retour / Terugbetaling"
#    now has 2 entries and moves above "Bestelling / Levering" (still 1)
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A5").Value = "Retour / Terugbetaling"
$dashboard.Range("B5").Value = 2
$dashboard.Range("A6").Value = "Bestelling / Levering"
$dashboard.Range("B6").Value = 1
